$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the row that held the 369/396 interval (old row 13) - it was
#    dropped from the data set.
$ws.Rows("13:13").Delete()

# 2. Insert a brand-new column A; everything that used to live in
#    A:E now shifts right into B:F, making room for a running index
#    column (used to drive the new "total heatmap" generation).
$ws.Columns("A:A").Insert()

# 3. Give the new index column (A2:A14) the same bold / bordered /
#    centered look that the header row (now B1:F1) already uses, by
#    copying B1's format onto it.
$ws.Range("B1").Copy()
$ws.Range("A2:A14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 4. Fill the new index column with sequential ids 0-12 (one per
#    remaining data row).
$lastRow = 14
$idx = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $idx
    $idx = $idx + 1
}

# 5. Fix the "Recusrive Step" typo -> "Iteration Step" everywhere it
#    appears (column F, after the column insert shifted E -> F).
$ws.Range("F10").Value = "Iteration Step"
$ws.Range("F11").Value = "Iteration Step"
$ws.Range("F12").Value = "Iteration Step"
